$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Departure Date for the sample row from "21 January 2025" to "8 February 2025"
$ws.Range("C2").Value = "8 February 2025"

# Reflect the active cell selection that was recorded in the saved file
$ws.Activate()
$ws.Range("I6").Select()

